# Refresh leve-profit calculation columns (currentAveragePrice* / LevePrice* /
# LeveProfit*, columns H:N) on each crafting-job sheet with the latest figures
# pulled by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")

# row 9 ("Distill, My Heart")
$ws.Range("H9").Value = 5350.5713
$ws.Range("I9").Value = 6250.706
$ws.Range("J9").Value = 1525
$ws.Range("K9").Value = 6250.706
$ws.Range("L9").Value = 1525
$ws.Range("M9").Value = -6081.706
$ws.Range("N9").Value = -1863

# row 17 ("One for the Road")
$ws.Range("H17").Value = 2402544.2
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 2491508.8
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 7474526.399999999
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -7474862.399999999

# row 33 ("Glazed and Confused")
$ws.Range("H33").Value = 3000
$ws.Range("I33").Value = 2666.6667
$ws.Range("K33").Value = 2666.6667
$ws.Range("M33").Value = -2437.6667

# row 64 ("Forged from the Void")
$ws.Range("H64").Value = 10249.125
$ws.Range("I64").Value = 10284.714
$ws.Range("K64").Value = 10284.714
$ws.Range("M64").Value = -10036.714

# row 67 ("Dodging the Draft (L)")
$ws.Range("H67").Value = 10249.125
$ws.Range("I67").Value = 10284.714
$ws.Range("K67").Value = 10284.714
$ws.Range("M67").Value = -9426.714

# row 100 ("Asking for a Friend")
$ws.Range("H100").Value = 3443.4285
$ws.Range("I100").Value = 3101.6667
$ws.Range("J100").Value = 3699.75
$ws.Range("K100").Value = 3101.6667
$ws.Range("L100").Value = 3699.75
$ws.Range("M100").Value = -2560.6667
$ws.Range("N100").Value = -4781.75

# row 107 ("Another Man's Ink")
$ws.Range("H107").Value = 919.4375
$ws.Range("I107").Value = 622.2143
$ws.Range("K107").Value = 622.2143
$ws.Range("M107").Value = 1297.7857

# row 116 ("Growing Up")
$ws.Range("H116").Value = 5432.4
$ws.Range("I116").Value = 4699.273
$ws.Range("K116").Value = 4699.273
$ws.Range("M116").Value = -1257.273

# row 132 ("Fast-forwarding Flora")
$ws.Range("H132").Value = 3915.037
$ws.Range("I132").Value = 3880.24
$ws.Range("K132").Value = 11640.72
$ws.Range("M132").Value = -9110.719999999999


# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")

# row 32 ("Ingot We Trust")
$ws.Range("H32").Value = 10529755
$ws.Range("I32").Value = 12503146
$ws.Range("K32").Value = 12503146
$ws.Range("M32").Value = -12502859

# row 43 ("They've Got Legs")
$ws.Range("H43").Value = 250034500
$ws.Range("I43").Value = 500020000
$ws.Range("J43").Value = 48999.5
$ws.Range("K43").Value = 500020000
$ws.Range("L43").Value = 48999.5
$ws.Range("M43").Value = -500019687
$ws.Range("N43").Value = -49625.5

# row 45 ("Hollow Hallmarks")
$ws.Range("H45").Value = 2615.9
$ws.Range("I45").Value = 2610
$ws.Range("J45").Value = 2629.6667
$ws.Range("K45").Value = 2610
$ws.Range("L45").Value = 2629.6667
$ws.Range("M45").Value = -2233
$ws.Range("N45").Value = -3383.6667

# row 61 ("Dealing with the Tough Stuff")
$ws.Range("H61").Value = 35718544
$ws.Range("I61").Value = 47622040
$ws.Range("K61").Value = 47622040
$ws.Range("M61").Value = -47621828

# row 74 ("As the Bolt Flies")
$ws.Range("H74").Value = 40047850
$ws.Range("I74").Value = 43529984
$ws.Range("K74").Value = 43529984
$ws.Range("M74").Value = -43529110

# row 77 ("Heavy Metal Banned (L)")
$ws.Range("H77").Value = 40047850
$ws.Range("I77").Value = 43529984
$ws.Range("K77").Value = 217649920
$ws.Range("M77").Value = -217645552

# row 97 ("Ore for Me")
$ws.Range("H97").Value = 1043.5
$ws.Range("I97").Value = 1116.2307
$ws.Range("J97").Value = 728.3333
$ws.Range("K97").Value = 1116.2307
$ws.Range("L97").Value = 728.3333
$ws.Range("M97").Value = -620.2307000000001
$ws.Range("N97").Value = -1720.3333

# row 136 ("Metal with Mettle")
$ws.Range("H136").Value = 35718544
$ws.Range("I136").Value = 47622040
$ws.Range("K136").Value = 142866120
$ws.Range("M136").Value = -142863570


# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")

# row 99 ("Meddle in Metal")
$ws.Range("H99").Value = 4470.2144
$ws.Range("I99").Value = 2890.25
$ws.Range("K99").Value = 2890.25
$ws.Range("M99").Value = -1392.25

# row 105 ("Ingot to Wing It")
$ws.Range("H105").Value = 10686.75
$ws.Range("I105").Value = 10686.75
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 10686.75
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -8939.75
$ws.Range("N105").ClearContents()


# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")

# row 58 ("You Do the Heavy Lifting")
$ws.Range("H58").Value = 2849.3547
$ws.Range("I58").Value = 2157.9092
$ws.Range("J58").Value = 4539.5557
$ws.Range("K58").Value = 2157.9092
$ws.Range("L58").Value = 4539.5557
$ws.Range("M58").Value = -1954.9092
$ws.Range("N58").Value = -4945.5557

# row 134 ("Wood You Be Quiet")
$ws.Range("H134").Value = 4615.077
$ws.Range("I134").Value = 4391.0967
$ws.Range("J134").Value = 5483
$ws.Range("K134").Value = 13173.2901
$ws.Range("L134").Value = 16449
$ws.Range("M134").Value = -10638.2901
$ws.Range("N134").Value = -21519

# row 136 ("Turali Quality")
$ws.Range("H136").Value = 2849.3547
$ws.Range("I136").Value = 2157.9092
$ws.Range("J136").Value = 4539.5557
$ws.Range("K136").Value = 6473.7276
$ws.Range("L136").Value = 13618.6671
$ws.Range("M136").Value = -3923.7276
$ws.Range("N136").Value = -18718.6671

# row 141 ("No Greater Treasure")
$ws.Range("H141").Value = 105613.6
$ws.Range("J141").Value = 123267
$ws.Range("L141").Value = 123267
$ws.Range("N141").Value = -133627


# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")

# row 26 ("A Grape Idea")
$ws.Range("H26").Value = 266
$ws.Range("I26").Value = 197.25
$ws.Range("J26").Value = 291
$ws.Range("K26").Value = 591.75
$ws.Range("L26").Value = 873
$ws.Range("M26").Value = -303.75
$ws.Range("N26").Value = -1449

# row 121 ("A Cookie for Your Troubles")
$ws.Range("H121").Value = 365.8
$ws.Range("I121").Value = 157.25
$ws.Range("J121").Value = 1200
$ws.Range("K121").Value = 471.75
$ws.Range("L121").Value = 3600
$ws.Range("M121").Value = 838.25
$ws.Range("N121").Value = -6220

# row 129 ("Comfort Food")
$ws.Range("H129").Value = 1406.174
$ws.Range("I129").Value = 562.6667
$ws.Range("K129").Value = 1688.0001
$ws.Range("M129").Value = 3311.9999


# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")

# row 2 ("Copper and Robbers")
$ws.Range("H2").Value = 207.87805
$ws.Range("I2").Value = 67.68000000000001
$ws.Range("J2").Value = 426.9375
$ws.Range("K2").Value = 67.68000000000001
$ws.Range("L2").Value = 426.9375
$ws.Range("M2").Value = 45.31999999999999
$ws.Range("N2").Value = -652.9375

# row 11 ("A Ringing Success")
$ws.Range("H11").Value = 35890830
$ws.Range("J11").Value = 3000
$ws.Range("L11").Value = 3000
$ws.Range("N11").Value = -3278

# row 12 ("Horn of Plenty")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

# row 97 ("If I'd a Koppranickel for Every Time...")
$ws.Range("H97").Value = 2113.375
$ws.Range("I97").Value = 1939.4166
$ws.Range("K97").Value = 1939.4166
$ws.Range("M97").Value = -1443.4166

# row 102 ("Put the Metal to the Peddle")
$ws.Range("H102").Value = 2207.818
$ws.Range("I102").Value = 2089.25
$ws.Range("J102").Value = 2524
$ws.Range("K102").Value = 2089.25
$ws.Range("L102").Value = 2524
$ws.Range("M102").Value = -467.25
$ws.Range("N102").Value = -5768

# row 109 ("You're My Wonderhall")
$ws.Range("H109").Value = 32997.6
$ws.Range("I109").Value = 32249.5
$ws.Range("J109").Value = 33496.332
$ws.Range("K109").Value = 32249.5
$ws.Range("L109").Value = 33496.332
$ws.Range("M109").Value = -31209.5
$ws.Range("N109").Value = -35576.332

# row 126 ("Gold Rush Order")
$ws.Range("H126").Value = 27368520
$ws.Range("I126").Value = 10105371
$ws.Range("J126").Value = 200000000
$ws.Range("K126").Value = 30316113
$ws.Range("L126").Value = 600000000
$ws.Range("M126").Value = -30313643
$ws.Range("N126").Value = -600004940

# row 132 ("On Board for Lar")
$ws.Range("H132").Value = 1517.3334
$ws.Range("I132").Value = 1456.5
$ws.Range("K132").Value = 4369.5
$ws.Range("M132").Value = -1839.5


# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")

# row 40 ("Best Served Toad")
$ws.Range("H40").Value = 6699.143
$ws.Range("I40").Value = 5816.5
$ws.Range("J40").Value = 11995
$ws.Range("K40").Value = 5816.5
$ws.Range("L40").Value = 11995
$ws.Range("M40").Value = -5680.5
$ws.Range("N40").Value = -12267

# row 69 ("Maybe He's a Lion")
$ws.Range("H69").Value = 36999
$ws.Range("I69").Value = 36999
$ws.Range("K69").Value = 36999
$ws.Range("M69").Value = -36188

# row 72 ("The Wyvern of It (L)")
$ws.Range("H72").Value = 36999
$ws.Range("I72").Value = 36999
$ws.Range("K72").Value = 110997
$ws.Range("M72").Value = -106941


# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")

# row 96 ("Skills on Display")
$ws.Range("H96").Value = 4793.25
$ws.Range("I96").Value = 2870.1
$ws.Range("K96").Value = 2870.1
$ws.Range("M96").Value = -1497.1

# row 126 ("A Polished Purchase")
$ws.Range("H126").Value = 9359.333000000001
$ws.Range("I126").Value = 11362
$ws.Range("J126").Value = 2350
$ws.Range("K126").Value = 34086
$ws.Range("L126").Value = 7050
$ws.Range("M126").Value = -31616
$ws.Range("N126").Value = -11990
